$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.392.32"
$ws.Range("E2").Value = "  +1.33%  "

$ws.Range("D3").Value = "'1.858.74"
$ws.Range("E3").Value = "  +0.44%  "

$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").Value = "'311.14"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("D6").Value = "'1.011"
$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("D7").Value = "'0.4771"
$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("D8").Value = "'0.3803"
$ws.Range("E8").Value = "  +3.46%  "

$ws.Range("D9").Value = "'0.07306"
$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("D10").Value = "'0.9300"
$ws.Range("E10").Value = "  -0.28%  "

$ws.Range("D11").Value = "'20.72"
$ws.Range("E11").Value = "  +4.87%  "

$ws.Range("D12").Value = "'0.07802"
$ws.Range("E12").Value = "  +0.72%  "

$ws.Range("D13").Value = "'1.871.81"
$ws.Range("E13").Value = "  +1.57%  "

$ws.Range("D14").Value = "'5.433"
$ws.Range("E14").Value = "  +1.57%  "

$ws.Range("D15").Value = "'6.543"
$ws.Range("E15").Value = "  +1.30%  "

$ws.Range("D16").Value = "'90.09"
$ws.Range("E16").Value = "  +1.09%  "

$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").Value = "'0.000008797"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("D20").Value = "'27.470.21"
$ws.Range("E20").Value = "  +1.56%  "

$ws.Range("D21").Value = "'14.63"
$ws.Range("E21").Value = "  +0.79%  "

$ws.Range("D22").Value = "'5.098"
$ws.Range("E22").Value = "  +0.46%  "

$ws.Range("D23").Value = "'10.70"
$ws.Range("E23").Value = "  +0.39%  "

$ws.Range("D24").Value = "'1.946"
$ws.Range("E24").Value = "  +1.10%  "

$ws.Range("D25").Value = "'155.36"
$ws.Range("E25").Value = "  +1.39%  "

$ws.Range("D26").Value = "'18.45"
$ws.Range("E26").Value = "  +1.09%  "

$ws.Range("D27").Value = "'2.006"
$ws.Range("E27").Value = "  -0.51%  "

$ws.Range("D28").Value = "'115.26"
$ws.Range("E28").Value = "  +0.76%  "

$ws.Range("D29").Value = "'4.934"
$ws.Range("E29").Value = "  -0.79%  "

$ws.Range("D30").Value = "'0.08880"
$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("D31").Value = "'3.326"
$ws.Range("E31").Value = "  +0.55%  "

$ws.Range("D32").Value = "'1.206"
$ws.Range("E32").Value = "  +2.05%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.587"
$ws.Range("E33").Value = "  +1.78%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7499"
$ws.Range("E34").Value = "  +1.16%  "

$ws.Range("D35").Value = "'2.702"
$ws.Range("E35").Value = "  -2.30%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'1.125"
$ws.Range("E36").Value = "  +0.95%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02041"
$ws.Range("E37").Value = "  +4.12%  "

$ws.Range("D38").Value = "'0.5532"
$ws.Range("E38").Value = "  +4.98%  "

$ws.Range("D39").Value = "'0.05260"
$ws.Range("E39").Value = "  -0.25%  "

$ws.Range("D40").Value = "'2.986"
$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("D41").Value = "'7.039"
$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("D42").Value = "'8.612"
$ws.Range("E42").Value = "  +4.44%  "

$ws.Range("D43").Value = "'0.1522"
$ws.Range("E43").Value = "  +0.37%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'10.71"
$ws.Range("E44").Value = "  +0.93%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.4876"
$ws.Range("E45").Value = "  +2.52%  "

$ws.Range("D46").Value = "'1.012"
$ws.Range("E46").Value = "  -0.18%  "

$ws.Range("D47").Value = "'1.660"
$ws.Range("E47").Value = "  +2.92%  "

$ws.Range("E48").Value = "  +1.02%  "

$ws.Range("D49").Value = "'67.42"
$ws.Range("E49").Value = "  +2.47%  "

$ws.Range("D50").Value = "'0.06079"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("D51").Value = "'0.9118"
$ws.Range("E51").Value = "  +2.44%  "

